$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(2,3,4,5,6,9,10,11,12,13,14)

$data = @(
  @(1.02, 1.044227090051403, 1.046391740782104, 1.042091236798147, 1.046729152420632, 1.042330235191239, 1.049293306182361, 1.049157306662782, 1.044868904457935, 1.049493773375442, 1.050783423032765),
  @(1.02, 1.045429885780342, 1.04734340325694, 1.043124335188844, 1.049059953381884, 1.042747866469422, 1.050142017139511, 1.049920238418817, 1.045712180833857, 1.051632336226778, 1.051633339256823),
  @(1.02, 1.046207295623237, 1.047958293717829, 1.043792298178111, 1.050562202015509, 1.043016144008523, 1.050689794832782, 1.050412388574341, 1.046256715367234, 1.053009889743201, 1.052181894857248),
  @(1.02, 1.046533912456415, 1.048216581630115, 1.044072988119318, 1.051192362594497, 1.043128462006152, 1.050919749873429, 1.050618928286177, 1.04648537204913, 1.053587553287666, 1.05241217646043),
  @(1.02, 1.046588740887692, 1.048259936963423, 1.044120110144567, 1.05129808916912, 1.04314729345418, 1.050958341038499, 1.050653586179627, 1.046523749018524, 1.053684460900404, 1.052450822429381),
  @(1.02, 1.046211660697682, 1.047961745802171, 1.043796049244693, 1.050570627644511, 1.043017646633726, 1.050692868798009, 1.050415149778985, 1.046259771728579, 1.05301761420295, 1.052184973187859),
  @(1.02, 1.044633765173786, 1.046713547574071, 1.042440486860529, 1.047518108430264, 1.042471782730781, 1.049580422862132, 1.049415459168707, 1.045154127660173, 1.050217823162451, 1.051070947451153),
  @(1.02, 1.041846382980668, 1.044507036419624, 1.040047689393825, 1.042092080690332, 1.041494777817012, 1.047609317860638, 1.047642110614726, 1.043197119379714, 1.045234912552862, 1.049097043254469),
  @(1.02, 1.039983195159899, 1.043031100895397, 1.038449525719662, 1.038440748977858, 1.040833090960958, 1.046287752385722, 1.046451764845028, 1.041886395196169, 1.041877649649973, 1.047773601004992),
  @(1.02, 1.039175181836479, 1.042390788923776, 1.037756755607229, 1.036851102038168, 1.040544076829087, 1.045713673772669, 1.04593436133139, 1.041317358771514, 1.040415053281613, 1.047198707134462),
  @(1.02, 1.038874857807975, 1.042152760780297, 1.037499312621093, 1.036259299727287, 1.040436344945706, 1.045500155885091, 1.045741873661718, 1.04110576650697, 1.039870403601307, 1.046984886026991),
  @(1.02, 1.038939287104717, 1.042203827189371, 1.037554540388417, 1.03638630455214, 1.040459471017984, 1.045545968915618, 1.045783176606288, 1.041151164091058, 1.039987295771566, 1.047030764117277),
  @(1.02, 1.039150360905458, 1.042371117297857, 1.037735477705993, 1.036802211006189, 1.04053517943431, 1.045696030044833, 1.045918456404194, 1.041299873137861, 1.040370060691634, 1.047181038350508),
  @(1.02, 1.039280384845491, 1.04247416524852, 1.037846943494046, 1.03705828621373, 1.040581775534335, 1.045788450489398, 1.046001766747454, 1.041391467521229, 1.040605711146408, 1.047273590042704),
  @(1.02, 1.040036793734937, 1.04307357018995, 1.038495486355332, 1.038546063659732, 1.040852218906912, 1.046325813129232, 1.046486061232569, 1.04192412867485, 1.041974526687161, 1.047811715799124),
  @(1.02, 1.04051093301444, 1.043449230959015, 1.038902095345041, 1.039476974052604, 1.041021189171301, 1.046662393021187, 1.046789314475613, 1.042257852905367, 1.042830741663812, 1.048148773673132),
  @(1.02, 1.040787371101018, 1.043668229895262, 1.039139190942457, 1.040019132015425, 1.041119505767252, 1.046858537906, 1.046966006593473, 1.042452365660362, 1.043329303532319, 1.048345197106179),
  @(1.02, 1.040881609268708, 1.043742883044445, 1.039220022226609, 1.040203855067459, 1.041152988434698, 1.046925388438598, 1.047026221872688, 1.042518665351833, 1.043499156686584, 1.048412142574202),
  @(1.02, 1.040460074740017, 1.04340893831127, 1.03885847758095, 1.039377181972218, 1.041003085204974, 1.046626299468447, 1.046756797992947, 1.0422220622666, 1.042738966520308, 1.048112628863407),
  @(1.02, 1.039088210274041, 1.042321859747504, 1.037682199446267, 1.036679774186714, 1.040512895689281, 1.0456518485336, 1.045878628216936, 1.041256088321785, 1.040257384297462, 1.047136794096461),
  @(1.02, 1.038224551562139, 1.041637281835917, 1.036941945918457, 1.034976050398721, 1.040202498484389, 1.045037553109719, 1.045324744638382, 1.040647428033963, 1.038689130934612, 1.046521626302618),
  @(1.02, 1.038682500748249, 1.042000294118558, 1.037334434179749, 1.035879977127927, 1.040367255272969, 1.045363357689439, 1.045618535443388, 1.040970216222117, 1.039521262816038, 1.046847893562211),
  @(1.02, 1.040483055774103, 1.043427145188163, 1.03887818679552, 1.039422276269282, 1.041011266352495, 1.04664260913681, 1.046771491380226, 1.042238234957595, 1.042780438393847, 1.048128961693369),
  @(1.02, 1.042567837222401, 1.045078325888216, 1.040666793214606, 1.043500656195772, 1.041749167751961, 1.048120201144874, 1.048101978709529, 1.043704105692111, 1.046529179552545, 1.049608652051567)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $excelRow = $startRow + $i
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($excelRow, $cols[$j]).Value = $row[$j]
    }
}
